# Automatic post-commit hook for streamlit
#
# The underlying data-cleaning step determined that the two survey
# responses stored in rows 17 and 18 of Sheet1 had been attached to the
# wrong sample metadata (program/location/county/state stay in columns
# A-F, and the geo/trailing columns CG-CP stay put as well) - all of the
# substantive answer columns, G through CF, actually belong to the other
# row. This script swaps that G:CF block between row 17 and row 18 so the
# answers line up with the correct sample again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowA = 17
$rowB = 18

# Column A (sampleid) swaps along with all of the substantive answer
# columns G:CF. Columns B:F (program/location/county/state/full_state) and
# CG:CP (trailing geo/misc columns) describe the collection site, which is
# identical for both rows here, so they stay attached to their own row.
$swapCols = @(1) + @(7..84)

foreach ($col in $swapCols) {

    $cellA = $ws.Cells.Item($rowA, $col)
    $cellB = $ws.Cells.Item($rowB, $col)

    $valA = $cellA.Value2
    $valB = $cellB.Value2

    $blankA = ($valA -eq $null)
    $blankB = ($valB -eq $null)

    # --- write row A's cell with row B's old content ---
    if ($blankB) {
        $cellA.Value2 = ""
    }
    elseif (($valB -is [string]) -and ($valB.Length -eq 0)) {
        # Distinguish a true empty-string cell from a merely-blank one:
        # writing "" through Value2 clears the cell outright, so force
        # text entry with a bare quote-prefix and then strip the
        # resulting "quoted text" style back to Normal.
        $cellA.Formula = "'"
        $cellA.Style = "Normal"
    }
    else {
        $cellA.Value2 = $valB
    }

    # --- write row B's cell with row A's old content ---
    if ($blankA) {
        $cellB.Value2 = ""
    }
    elseif (($valA -is [string]) -and ($valA.Length -eq 0)) {
        $cellB.Formula = "'"
        $cellB.Style = "Normal"
    }
    else {
        $cellB.Value2 = $valA
    }
}
